{"js": "const replacements = [\n  [\"27\u00d757=1539\", \"71\u00d759=4189\"],\n  [\"94\u00d732=3008\", \"92\u00d757=5244\"],\n  [\"19\u00d732=608\", \"98\u00d783=8134\"],\n  [\"26\u00d730=780\", \"44\u00d770=3080\"],\n  [\"37\u00d743=1591\", \"74\u00d772=5328\"],\n  [\"88\u00d730=2640\", \"73\u00d789=6497\"],\n  [\"73\u00d770=5110\", \"70\u00d716=1120\"],\n  [\"43\u00d766=2838\", \"41\u00d730=1230\"],\n  [\"41\u00d754=2214\", \"27\u00d729=783\"],\n  [\"78\u00d735=2730\", \"86\u00d726=2236\"],\n  [\"81\u00d782=6642\", \"82\u00d772=5904\"],\n  [\"12\u00d786=1032\", \"64\u00d765=4160\"],\n  [\"91\u00d771=6461\", \"85\u00d767=5695\"],\n  [\"74\u00d793=6882\", \"88\u00d789=7832\"],\n  [\"48\u00d751=2448\", \"49\u00d747=2303\"],\n  [\"92\u00d793=8556\", \"94\u00d761=5734\"],\n  [\"18\u00d754=972\", \"86\u00d740=3440\"],\n  [\"86\u00d741=3526\", \"84\u00d741=3444\"],\n  [\"16\u00d713=208\", \"74\u00d786=6364\"],\n  [\"20\u00d730=600\", \"60\u00d772=4320\"],\n  [\"91\u00d749=4459\", \"40\u00d739=1560\"],\n  [\"34\u00d726=884\", \"52\u00d717=884\"],\n  [\"99\u00d777=7623\", \"77\u00d790=6930\"],\n  [\"74\u00d721=1554\", \"33\u00d726=858\"],\n  [\"70\u00d718=1260\", \"96\u00d791=8736\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"27\u00d757=1539\", \"71\u00d759=4189\"),\n    @(\"94\u00d732=3008\", \"92\u00d757=5244\"),\n    @(\"19\u00d732=608\", \"98\u00d783=8134\"),\n    @(\"26\u00d730=780\", \"44\u00d770=3080\"),\n    @(\"37\u00d743=1591\", \"74\u00d772=5328\"),\n    @(\"88\u00d730=2640\", \"73\u00d789=6497\"),\n    @(\"73\u00d770=5110\", \"70\u00d716=1120\"),\n    @(\"43\u00d766=2838\", \"41\u00d730=1230\"),\n    @(\"41\u00d754=2214\", \"27\u00d729=783\"),\n    @(\"78\u00d735=2730\", \"86\u00d726=2236\"),\n    @(\"81\u00d782=6642\", \"82\u00d772=5904\"),\n    @(\"12\u00d786=1032\", \"64\u00d765=4160\"),\n    @(\"91\u00d771=6461\", \"85\u00d767=5695\"),\n    @(\"74\u00d793=6882\", \"88\u00d789=7832\"),\n    @(\"48\u00d751=2448\", \"49\u00d747=2303\"),\n    @(\"92\u00d793=8556\", \"94\u00d761=5734\"),\n    @(\"18\u00d754=972\", \"86\u00d740=3440\"),\n    @(\"86\u00d741=3526\", \"84\u00d741=3444\"),\n    @(\"16\u00d713=208\", \"74\u00d786=6364\"),\n    @(\"20\u00d730=600\", \"60\u00d772=4320\"),\n    @(\"91\u00d749=4459\", \"40\u00d739=1560\"),\n    @(\"34\u00d726=884\", \"52\u00d717=884\"),\n    @(\"99\u00d777=7623\", \"77\u00d790=6930\"),\n    @(\"74\u00d721=1554\", \"33\u00d726=858\"),\n    @(\"70\u00d718=1260\", \"96\u00d791=8736\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
